$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("DREXvsHFD", "ImmSC", 34, 20, 14, 58.8, 183, 70, 113, 38.3),
    @("DREXvsHFD", "mySC", 59, 34, 25, 57.6, 227, 94, 133, 41.4),
    @("DREXvsHFD", "nmSC", 44, 29, 15, 65.9, 78, 35, 43, 44.9),
    @("DRvsHFD", "ImmSC", 15, 4, 11, 26.7, 45, 21, 24, 46.7),
    @("DRvsHFD", "mySC", 75, 19, 56, 25.3, 301, 201, 100, 66.8),
    @("DRvsHFD", "nmSC", 6, 1, 5, 16.7, 29, 20, 9, 69),
    @("EXvsHFD", "ImmSC", 25, 11, 14, 44, 105, 78, 27, 74.3),
    @("EXvsHFD", "mySC", 101, 31, 70, 30.7, 517, 348, 169, 67.3),
    @("EXvsHFD", "nmSC", 39, 9, 30, 23.1, 106, 83, 23, 78.3),
    @("HFDvsSD", "ImmSC", 22, 14, 8, 63.6, 51, 17, 34, 33.3),
    @("HFDvsSD", "mySC", 52, 36, 16, 69.2, 195, 86, 109, 44.1),
    @("HFDvsSD", "nmSC", 27, 10, 17, 37, 81, 55, 26, 67.9),
    @("DREXvsHFD", "majorSC", 95, 58, 37, 61.1, 416, 145, 271, 34.9),
    @("DRvsHFD", "majorSC", 80, 25, 55, 31.2, 324, 204, 120, 63),
    @("EXvsHFD", "majorSC", 179, 58, 121, 32.4, 929, 585, 344, 63),
    @("HFDvsSD", "majorSC", 75, 44, 31, 58.7, 175, 108, 67, 61.7),
    @("DREXvsHFD", "aggSC", 106, 66, 40, 62.3, 403, 169, 234, 41.9),
    @("DRvsHFD", "aggSC", 90, 23, 67, 25.6, 363, 232, 131, 63.9),
    @("EXvsHFD", "aggSC", 134, 40, 94, 29.9, 636, 434, 202, 68.2),
    @("HFDvsSD", "aggSC", 87, 50, 37, 57.5, 307, 148, 159, 48.2)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $rowNum = $r + 2
    $rowData = $data[$r]
    for ($c = 0; $c -lt $rowData.Length; $c++) {
        $colNum = $c + 1
        $ws.Cells.Item($rowNum, $colNum).Value = $rowData[$c]
    }
}
